$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text before writing so that values
# such as "0.9999" or "26.914.43" are stored verbatim as strings instead
# of being auto-converted to numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.914.43"
$ws.Range("E2").Value = "  +4.42%  "
$ws.Range("D3").Value = "1.878.09"
$ws.Range("E3").Value = "  +3.53%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "278.95"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "0.5278"
$ws.Range("E7").Value = "  +4.24%  "
$ws.Range("D8").Value = "0.3450"
$ws.Range("E8").Value = "  -1.66%  "
$ws.Range("D9").Value = "45.10"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").Value = "0.06959"
$ws.Range("E10").Value = "  +4.31%  "
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("D12").Value = "0.8063"
$ws.Range("E12").Value = "  -2.93%  "
$ws.Range("D13").Value = "0.07855"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "1.856.90"
$ws.Range("E14").Value = "  +2.41%  "
$ws.Range("D15").Value = "5.180"
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").Value = "90.18"
$ws.Range("E16").Value = "  +3.03%  "
$ws.Range("D17").Value = "14.61"
$ws.Range("E17").Value = "  +3.98%  "
$ws.Range("D18").Value = "0.9993"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "0.000008102"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").Value = "0.9996"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").Value = "26.957.24"
$ws.Range("E21").Value = "  +4.43%  "
$ws.Range("D22").Value = "2.151.11"
$ws.Range("E22").Value = "  +5.39%  "
$ws.Range("D23").Value = "4.749"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").Value = "6.184"
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("D26").Value = "2.365"
$ws.Range("E26").Value = "  +9.33%  "
$ws.Range("D27").Value = "145.70"
$ws.Range("E27").Value = "  +2.62%  "
$ws.Range("E28").Value = "  +2.02%  "
$ws.Range("D29").Value = "1.662"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").Value = "113.82"
$ws.Range("E30").Value = "  +3.86%  "
$ws.Range("D31").Value = "4.383"
$ws.Range("E31").Value = "  +1.30%  "
$ws.Range("D32").Value = "4.335"
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("D33").Value = "0.08906"
$ws.Range("E33").Value = "  +1.08%  "
$ws.Range("D34").Value = "0.04948"
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("D35").Value = "1.179"
$ws.Range("E35").Value = "  +3.94%  "
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("D38").Value = "3.274"
$ws.Range("E38").Value = "  +3.86%  "
$ws.Range("E39").Value = "  +4.43%  "
$ws.Range("D40").Value = "0.01853"
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("D41").Value = "0.5166"
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("D42").Value = "0.9584"
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("D43").Value = "116.42"
$ws.Range("E43").Value = "  +2.73%  "
$ws.Range("D44").Value = "6.211"
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("D45").Value = "8.145"
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("D46").Value = "0.9991"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").Value = "0.4507"
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("D48").Value = "0.1347"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").Value = "9.371"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").Value = "36.27"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").Value = "0.05954"
$ws.Range("E51").Value = "  +2.09%  "

# Restore the original (unstyled) formatting of the data cells.
$dataRange.ClearFormats()
